# FAKER API to generate unique random test data Implemented
$wb = $excel.ActiveWorkbook

# The fake-data generator added a unique value ("aaaaaa") into the
# fName column (A2) of Sheet2, which was previously left blank.
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("A2").Value = "aaaaaa"

# Make Sheet2 the active / selected sheet (it was Sheet4 before).
$ws2.Activate()
$ws2.Select()
